$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.598.27"
$ws.Range("E2").Value = "  -3.21%  "
$ws.Range("D3").Value = "1.847.67"
$ws.Range("E3").Value = "  -4.11%  "
$ws.Range("E4").Value = "  -0.99%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4659"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3902"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.26"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07899"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9797"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "1.893.49"
$ws.Range("E13").Value = "  -2.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.835"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.989"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06902"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.35%  "
$ws.Range("E19").Value = "  -3.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("D22").Value = "28.617.13"
$ws.Range("E22").Value = "  -3.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.389"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.03%  "
$ws.Range("E24").Value = "  -6.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.162"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("D26").Value = "2.101.56"
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("E28").Value = "  -2.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.089"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.027"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9738"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09346"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.355"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.21%  "
$ws.Range("E35").Value = "  -2.19%  "
$ws.Range("E36").Value = "  -2.82%  "
$ws.Range("E37").Value = "  -2.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02198"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.167"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5707"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.665"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1794"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.369"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.249"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5382"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.07109"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.903"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "42.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.27%  "
